$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New comment for row 2 (this introduces shared string index 16)
$ws.Range("G2").Value = "Constant duty cycle input, used for system ID"

# Fill the ditto marks down column G for rows 3-7 (reuses shared string index 8 "^")
$ws.Range("G3").Value = "^"
$ws.Range("G4").Value = "^"
$ws.Range("G5").Value = "^"
$ws.Range("G6").Value = "^"
$ws.Range("G7").Value = "^"

# New row 9 - OpenLoopSimulinkData_ID08
$ws.Range("A9").Value = "OpenLoopSimulinkData_ID08"
$ws.Range("B9").Value = "^"
$ws.Range("C9").Value = 900
$ws.Range("D9").Value = "^"

# New row 10 - OpenLoopSimulinkData_ID09 (introduced before E9/G9 strings per original authoring order)
$ws.Range("A10").Value = "OpenLoopSimulinkData_ID09"

$ws.Range("E9").Value = "various"

$ws.Range("B10").Value = "^"
$ws.Range("C10").Value = "^"
$ws.Range("D10").Value = "^"
$ws.Range("E10").Value = "^"

$ws.Range("F9").Value = "^"
$ws.Range("F10").Value = "^"

$ws.Range("G10").Value = "with fan on low but close to heater"
$ws.Range("G9").Value = "open loop run (in box)"

# Update selection to reflect the author's final cursor position
$ws.Range("A10").Select()
